$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously empty row 3 (autogluon) results
$ws.Range("B3").Value = "0.101 (0.054 ± 0.025)"
$ws.Range("C3").Value = "00:03:57 (00:04:38 ± 00:00:18)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "'61"

# Fix mojibake "Â±" -> "±" in row 4 (autokeras)
$ws.Range("B4").Value = "0.332 (0.253 ± 0.033)"
$ws.Range("C4").Value = "00:00:45 (00:00:55 ± 00:00:06)"
$ws.Range("D4").Value = "00:00:03 (00:00:03 ± 00:00:00)"

# Fix mojibake "Â±" -> "±" in row 6 (autosklearn)
$ws.Range("B6").Value = "0.679 (0.565 ± 0.068)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 ± 00:00:03)"
$ws.Range("D6").Value = "00:00:00 (00:00:06 ± 00:00:04)"
